$d = $word.ActiveDocument

# =====================================================================
# Change 1: insert "Team Number:" / "Team Member Names:" paragraphs
#           right after "Keep answers as short as possible..." (both
#           styled BodyText, like the rest of the submission-sheet
#           boilerplate).
# =====================================================================
$introPara = $d.Paragraphs.Item(2)
if ($introPara.Range.Text -notlike "*Keep answers as short as possible*") {
    throw "unexpected paragraph 2 content: $($introPara.Range.Text)"
}

$introPara.Range.InsertParagraphAfter()
$teamNumberPara = $d.Paragraphs.Item(3)
$teamNumberPara.Range.Text = "Team Number:"
$teamNumberPara.Style = "BodyText"

$teamNumberPara.Range.InsertParagraphAfter()
$teamMembersPara = $d.Paragraphs.Item(4)
$teamMembersPara.Range.Text = "Team Member Names:"
$teamMembersPara.Style = "BodyText"

# =====================================================================
# Change 2: append a new Compact checklist item ("Results are used to
#           explain measurement errors...") right after the *second*
#           "Calculations are correct and yield the proper Teensy unit
#           value." checklist entry (the one belonging to the complete
#           specification of the 175 kHz / 200 Hz measurement section),
#           continuing the same numbered list.
# =====================================================================
$paras = $d.Paragraphs
$matchCount = 0
$calcPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "*Calculations are correct and yield the proper Teensy unit value.*") {
        $matchCount = $matchCount + 1
        if ($matchCount -eq 2) {
            $calcPara = $cand
        }
    }
}
if ($calcPara -eq $null) {
    throw "could not find the second 'Calculations are correct' checklist paragraph"
}

$calcListTemplate = $calcPara.Range.ListFormat.ListTemplate
$calcPara.Range.InsertParagraphAfter()
$resultsPara = $calcPara.Next()
$resultsPara.Range.Text = "☐ Results are used to explain measurement errors that are present in the 175 kHz measurement that are not present in the 200 Hz measurement."
$resultsPara.Style = "Compact"
$resultsPara.Range.ListFormat.ApplyListTemplate($calcListTemplate, $true)

# =====================================================================
# Change 3: extend the extra-credit mini-report sentence with the
#           page-break instruction.
# =====================================================================
$replaced = $d.Content.Find.Execute(
    "of <1 page for each extra credit section that you attempt.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "of <1 page for each extra credit section that you attempt. Put a page break between each extra credit mini-report.",
    2)
if (-not $replaced) {
    throw "could not find the extra-credit mini-report sentence to extend"
}
